# Apply updated forecast-error values to existing rows (Q0..Q8, rows 2..10)
# and append a new row for Q9 (row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values: row index => B,C,D,E,F,G
$data = @{
    2  = @(0.005369470315216797, 1.076795149390729, 5.444793568182083, 2.333408144363537, 2.356620468321115, 51)
    3  = @(0.08004204401876436, 1.074134375449918, 5.376472252851192, 2.318722116350123, 2.340867082514619, 50)
    4  = @(0.01808692518996532, 1.091555239787779, 5.397618064264203, 2.323277440226243, 2.347282362305163, 49)
    5  = @(0.0673762441061526, 1.105808622426705, 5.589789055284754, 2.364273473032414, 2.388322547876662, 48)
    6  = @(0.02957733957197041, 1.131070466074804, 5.550086184435333, 2.355862089434637, 2.381143925742878, 47)
    7  = @(0.08923231662305202, 1.105571907883932, 5.573846403776829, 2.360899490401239, 2.385282018651314, 46)
    8  = @(0.04225835016865745, 1.172851906978339, 5.774705933445779, 2.403061783110409, 2.429840096861529, 45)
    9  = @(0.1053273583819827, 1.115688490554983, 5.76308808600117, 2.400643265044011, 2.426058849520155, 44)
    10 = @(-0.01140659214618309, 1.164923355123781, 5.955443118527865, 2.440377658996219, 2.469231899811152, 43)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}

# New row 11 for Q9
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Cells.Item(11, 1).Value = "Q9"

$ws.Cells.Item(11, 2).Value = -0.001181092721129312
$ws.Cells.Item(11, 3).Value = 1.078446272129183
$ws.Cells.Item(11, 4).Value = 5.792657664955174
$ws.Cells.Item(11, 5).Value = 2.406794063677899
$ws.Cells.Item(11, 6).Value = 2.435968096926322
$ws.Cells.Item(11, 7).Value = 42
